$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.240.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.306.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("E5").Value = '  +0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.20%  '

$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.663.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.352.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.079.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.72%  '

$ws.Range("E20").Value = '  +3.78%  '

$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.64%  '

$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.78%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '

$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0697'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("E38").Value = '  +1.42%  '

$ws.Range("E39").Value = '  +3.78%  '

$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.994.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("E44").Value = '  +4.37%  '

$ws.Range("E45").Value = '  +4.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.50%  '

$ws.Range("E49").Value = '  +5.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.530.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.00%  '

